$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 13, shifting rows 13-23 down to 14-24
$ws.Rows.Item(13).Insert()

# Remove the stray leftover cell at A13 (the inserted row has no label in column A)
$ws.Range("A13").Clear()

# Copy formatting (wrap-text / red-text styles) from row 11 (B/C) onto the new row 13 B/C cells
$ws.Range("B11:C11").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 10 (Objetivos:) - fix the value that used to (incorrectly) hold the teacher name
$ws.Range("B10").Value = 'Fornecer uma visão geral da indústria de serviços. Desenvolver projeto de criação de um novo serviço.'
$ws.Range("C10").Value = 'Fornecer uma visão geral da indústria de serviços. Desenvolver projeto de criação de um novo serviço.'

# Row 13 (new, no label) - teacher name now lives here, under Docentes responsaveis:
$ws.Range("B13").Value = '5840560 - Marco Antonio Carvalho Pereira'
$ws.Range("C13").Value = '5840560 - Marco Antonio Carvalho Pereira'

# Row 14 (Programa resumido:) - fix the value (used to say "Semestral")
$ws.Range("B14").Value = 'Introdução a Indústria de Serviços.Características Essenciais e diferenciadoras de Serviços.Projeto de Novo Serviço: Planejamento estratégico, Concepção do Serviço, Processos, Instalações. Avaliação e Melhoria.'
$ws.Range("C14").Value = 'Introdução a Indústria de Serviços.Características Essenciais e diferenciadoras de Serviços.Projeto de Novo Serviço: Planejamento estratégico, Concepção do Serviço, Processos, Instalações. Avaliação e Melhoria.'

# Row 16 (Programa:) - fix the value (used to say "01/01/2021")
$ws.Range("B16").Value = 'Características Essenciais e diferenciadoras de Serviços. Ciclo de Serviços.Projeto de Novo Serviço: Planejamento estratégico (Forças de Porter, Posicionamento Estratégico). Concepção do Serviço(Conceito de Serviço. Benchmarking. SERVQUAL. Geração e Seleção de Ideias. Pacote de Serviços. Especificações deServiço). Processos (Blue Print. Padronização. Entrega do Serviço. Recrutamento e Treinamento). Instalações (Seleção eLocalização. Gestão de Evidências Físicas. Projeto do Espaço Físico. Estudo da Capacidade Produtiva). Avaliação e Melhoria'
$ws.Range("C16").Value = 'Características Essenciais e diferenciadoras de Serviços. Ciclo de Serviços.Projeto de Novo Serviço: Planejamento estratégico (Forças de Porter, Posicionamento Estratégico). Concepção do Serviço(Conceito de Serviço. Benchmarking. SERVQUAL. Geração e Seleção de Ideias. Pacote de Serviços. Especificações deServiço). Processos (Blue Print. Padronização. Entrega do Serviço. Recrutamento e Treinamento). Instalações (Seleção eLocalização. Gestão de Evidências Físicas. Projeto do Espaço Físico. Estudo da Capacidade Produtiva). Avaliação e Melhoria'

# Row 19 (Metodo:) - fix the value (used to hold the teacher name)
$ws.Range("B19").Value = 'Aulas Expositivas; trabalhos em grupo; exercícios individuais e palestras.'
$ws.Range("C19").Value = 'Aulas Expositivas; trabalhos em grupo; exercícios individuais e palestras.'

# Row 20 (Criterio:) - fix the value (used to hold the Metodo text)
$ws.Range("B20").Value = 'Avaliação individual (Peso entre 20-40%) e do projeto realizado em equipe (peso entre 60-80%)'
$ws.Range("C20").Value = 'Avaliação individual (Peso entre 20-40%) e do projeto realizado em equipe (peso entre 60-80%)'

# Row 21 (Norma de recuperacao:) - fix the value (used to hold the Criterio text)
$ws.Range("B21").Value = 'NF = (MF + PR)/2, onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota darecuperação.'
$ws.Range("C21").Value = 'NF = (MF + PR)/2, onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota darecuperação.'

# Row 22 (Bibliografia:) - fix the value (used to hold the Norma de recuperacao text)
$ws.Range("B22").Value = 'CARVALHO, M. M. (organizadora) e outros. Gestão de Serviços: Casos Brasileiros. Atlas. 2013CORREA, H. C. e CAON, M. Gestão de Serviços: Lucratividade por meio de operação e de satisfação dos clientes. Atlas, 2014FITZSIMMONS, J.; FITZSIMMONS, M.J. Administração de serviços: operações, estratégia e tecnologia de informação. PortoAlegre: Bookman, 2000.GIANESI, I e CORREA, H. Administração Estratégia de Serviços, ATLAS, 1995 – SPJOHNSTON, R. e CLARK, G. Administração e Operações de Serviços. Atlas, 2001LOVELOCK, C.H.; WRIGHT, L. Serviços: marketing e gestão. São Paulo: Saraiva, 2001.MELLO, C. H. P.; NETO, P. L. O. C.; TURRIONI, J.B. SILVA, C. E. S. Gestão do Processo de Desenvolvimento de Serviços.Atlas. 2010NORMANN, R. Administração de Serviços. São Paulo. Atlas. 1992.Bibliografia complementar será indicada ao longo do curso.'
$ws.Range("C22").Value = 'CARVALHO, M. M. (organizadora) e outros. Gestão de Serviços: Casos Brasileiros. Atlas. 2013CORREA, H. C. e CAON, M. Gestão de Serviços: Lucratividade por meio de operação e de satisfação dos clientes. Atlas, 2014FITZSIMMONS, J.; FITZSIMMONS, M.J. Administração de serviços: operações, estratégia e tecnologia de informação. PortoAlegre: Bookman, 2000.GIANESI, I e CORREA, H. Administração Estratégia de Serviços, ATLAS, 1995 – SPJOHNSTON, R. e CLARK, G. Administração e Operações de Serviços. Atlas, 2001LOVELOCK, C.H.; WRIGHT, L. Serviços: marketing e gestão. São Paulo: Saraiva, 2001.MELLO, C. H. P.; NETO, P. L. O. C.; TURRIONI, J.B. SILVA, C. E. S. Gestão do Processo de Desenvolvimento de Serviços.Atlas. 2010NORMANN, R. Administração de Serviços. São Paulo. Atlas. 1992.Bibliografia complementar será indicada ao longo do curso.'

